$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [string][char]0x2083

$ws.Range("D2").Value = "40.273.12"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "2.235.85"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.73"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "89.15"
$ws.Range("E6").Value = "  +2.41%  "
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.483"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.55"
$ws.Range("E10").Value = "  -2.46%  "
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").Value = "2.584.29"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.92"
$ws.Range("E15").Value = "  -2.53%  "
$ws.Range("D16").Value = "2.237.90"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.737"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").Value = "40.214.76"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "0.0" + $sub3 + "0891"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.49"
$ws.Range("E20").Value = "  +7.42%  "
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.77"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.27"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.47"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E26").Value = "  -1.26%  "
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "155.48"
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.37"
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.92"
$ws.Range("E35").Value = "  +6.68%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.36"
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "15.92"
$ws.Range("E38").Value = "  -4.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0984"
$ws.Range("E39").Value = "  -2.79%  "
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").Value = "2.137.12"
$ws.Range("E41").Value = "  +5.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.86"
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.15"
$ws.Range("E43").Value = "  -3.57%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.23"
$ws.Range("E44").Value = "  +11.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0270"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.82"
$ws.Range("E46").Value = "  -1.44%  "
$ws.Range("E47").Value = "  +4.81%  "
$ws.Range("D48").Value = "2.448.99"
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("E49").Value = "  +1.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "69.66"
$ws.Range("E50").Value = "  -2.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "89.19"
$ws.Range("E51").Value = "  -0.75%  "
